$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Percent Started" header in column G
$ws.Range("G1").Value = "Percent Started"

# Percent Started values for rows 2..200 (player bench/starter usage %)
$percentStarted = @(
    95.58,95.89,8.17,15.93,42.59,98.36,6.43,38.23,74.62,51.58,61.47,26.17,30.26,1.22,43.75,26.03,15.08,98.97,42.75,47.36,
    75.58,23.08,46.18,34.89,21.61,27.72,2.02,94.13,12.65,12.29,62.69,32.52,26.68,17.63,99.31999999999999,86.56,37.6,36.5,
    84.01000000000001,73.42,41.17,7.65,9.369999999999999,24.87,36.12,57.76,18.18,72.56,24.1,30.53,12.2,69.31,56.63,0.28,
    74.03,42.78,26.06,43.17,43.01,48.05,78.38,9.619999999999999,3.67,53.84,7.69,40.95,4.43,64.36,16.43,83.01000000000001,
    52.64,72.66,72.13,49.4,28.32,91.73,37.19,41.6,70.70999999999999,53.34,52.43,19.79,4.64,16.01,71.98999999999999,31.34,
    89.47,0.84,22.09,0.08,24.59,94.02,8.02,55.49,17.87,27.8,5.19,63.11,1.08,1.26,17.19,44.51,72.34999999999999,55.95,4.91,
    14.81,73.41,16.88,7.43,64.79000000000001,18.9,46.18,29.68,12.77,66.97,14.56,37.99,75.62,94.44,79.42,24.1,28.52,14.37,
    65.84,58.82,46.37,68.52,10.76,13.5,7.51,20.73,6.97,30.39,73.87,61.65,4.21,82.40000000000001,17.35,0.51,84.72,
    0.6899999999999999,70.01000000000001,98.19,86.02,1.31,71.08,1.71,24.59,0.86,18.15,99.59,91.69,39.92,53.1,26.88,69.36,
    80.81999999999999,6.7,81.48999999999999,64.01000000000001,67.73999999999999,32.32,26.56,73.47,42.59,15.87,99.38,13.88,
    50.69,9.92,24.23,39.9,80.11,34.2,68.66,15.43,88.68000000000001,58.2,13.63,5.43,17.36,71.55,57.53,53.11,82.87,49.32,
    59.41,2.6,9.130000000000001,49.96,3.52,11.35,4.5,0.74,40.34,23.09,33.8,3.54,52.76
)

for ($i = 0; $i -lt $percentStarted.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $percentStarted[$i]
}
